$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Sheet1"

# 2. Insert a new column before C ("Thoi luong" shifts right, new "The loai" column is inserted)
$ws.Columns("C").Insert()

# 3. Set header text for new column C
$ws.Range("C3").Value = "Thể loại"

# 4. Resize columns B:I (A keeps its original width)
$ws.Columns("B").ColumnWidth = 29.714285714285715
$ws.Columns("C").ColumnWidth = 24.0
$ws.Columns("D").ColumnWidth = 11.714285714285714
$ws.Columns("E").ColumnWidth = 16.142857142857142
$ws.Columns("F").ColumnWidth = 22.0
$ws.Columns("G").ColumnWidth = 37.57142857142857
$ws.Columns("H").ColumnWidth = 80.0
$ws.Columns("I").ColumnWidth = 16.571428571428573

# 5. Vertically center the data row (row 4) cells
$ws.Range("A4:I4").VerticalAlignment = -4108

# 6. Reset scroll position and selection to match the saved view
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F12").Select()

Write-Output "done"
